# Rename the four logo pictures living in the document's headers/footers.
#
#   footer (first page)  - id=2 - PearsonLogo.png : image1.png -> image2.png
#   footer (default/odd) - id=4 - PearsonLogo.png : image1.png -> image2.png
#   header (first page)  - id=1 - BTec_Logo-Orange: image2.jpg -> image1.jpg
#   header (default/odd) - id=3 - BTec_Logo-Orange: image2.jpg -> image1.jpg
#
# InlineShape has no settable Name in the Word object model, so each
# picture is briefly converted to a floating Shape (which does expose
# Name), renamed, then converted back to an inline layout.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineLogo($shapeRangeOwner, $newName) {
    $inlineShape = $shapeRangeOwner.Range.InlineShapes.Item(1)
    $floating = $inlineShape.ConvertToShape()
    $floating.Name = $newName
    # wdWrapInline = 7 - put it back inline so the drawing stays a
    # <wp:inline> (not a floating <wp:anchor>).
    $floating.WrapFormat.Type = 7
}

# Footers carry the Pearson Edexcel logo (image1.png -> image2.png).
Rename-InlineLogo $sec.Footers.Item(2) "image2.png"   # footer1.xml (first page)
Rename-InlineLogo $sec.Footers.Item(1) "image2.png"   # footer2.xml (default)

# Headers carry the BTEC logo (image2.jpg -> image1.jpg).
Rename-InlineLogo $sec.Headers.Item(2) "image1.jpg"   # header1.xml (first page)
Rename-InlineLogo $sec.Headers.Item(1) "image1.jpg"   # header2.xml (default)
